# Update "想去人数" (people interested) counts on the 展览 and 全部类型 sheets
# to reflect newly fetched totals.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value for column F
$exhibitUpdates = @{
    3  = 1347
    5  = 101
    7  = 11634
    8  = 4373
    13 = 2543
    16 = 35
    17 = 5066
    18 = 61
    19 = 183
    21 = 11327
    22 = 11254
    24 = 45
    28 = 20
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$allUpdates = @{
    3  = 1347
    5  = 101
    7  = 11634
    8  = 4373
    13 = 2543
    17 = 35
    18 = 5066
    19 = 61
    20 = 183
    22 = 11327
    23 = 11254
    25 = 45
    29 = 20
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
